$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "per_si_allsi.ben_q1_tot"
$ws.Cells.Item(2, 2).Value = "SI.BEN.Q1.TOT"
$ws.Cells.Item(3, 1).Value = "per_sa_allsa.ben_q1_tot"
$ws.Cells.Item(3, 2).Value = "SA.BEN.Q1.TOT"
$ws.Cells.Item(4, 1).Value = "per_lm_alllm.ben_q1_tot"
$ws.Cells.Item(4, 2).Value = "LM.BEN.Q1.TOT"
$ws.Cells.Item(5, 1).Value = "per_allsp.ben_q1_tot"
$ws.Cells.Item(5, 2).Value = "SP.BEN.Q1.TOT"
$ws.Cells.Item(6, 1).Value = "per_allsp.adq_pop_tot"
$ws.Cells.Item(6, 2).Value = "SP.ADQ.POP.TOT"
$ws.Cells.Item(7, 1).Value = "per_si_allsi.adq_pop_tot"
$ws.Cells.Item(7, 2).Value = "SI.ADQ.POP.TOT"
$ws.Cells.Item(8, 1).Value = "per_lm_alllm.adq_pop_tot"
$ws.Cells.Item(8, 2).Value = "LM.ADQ.POP.TOT"
$ws.Cells.Item(9, 1).Value = "per_sa_allsa.adq_pop_tot"
$ws.Cells.Item(9, 2).Value = "SA.ADQ.POP.TOT"
$ws.Cells.Item(10, 1).Value = "per_sa_allsa.cov_q2_tot"
$ws.Cells.Item(10, 2).Value = "SA.COV.Q2.TOT"
$ws.Cells.Item(11, 1).Value = "per_lm_alllm.cov_q5_tot"
$ws.Cells.Item(11, 2).Value = "LM.COV.Q5.TOT"
$ws.Cells.Item(12, 1).Value = "per_lm_alllm.cov_q1_tot"
$ws.Cells.Item(12, 2).Value = "LM.COV.Q1.TOT"
$ws.Cells.Item(13, 1).Value = "per_lm_alllm.cov_q4_tot"
$ws.Cells.Item(13, 2).Value = "LM.COV.Q4.TOT"
$ws.Cells.Item(14, 1).Value = "per_sa_allsa.cov_q4_tot"
$ws.Cells.Item(14, 2).Value = "SA.COV.Q4.TOT"
$ws.Cells.Item(15, 1).Value = "per_si_allsi.cov_pop_tot"
$ws.Cells.Item(15, 2).Value = "SI.COV.POP.TOT"
$ws.Cells.Item(16, 1).Value = "per_si_allsi.cov_q2_tot"
$ws.Cells.Item(16, 2).Value = "SI.COV.Q2.TOT"
$ws.Cells.Item(17, 1).Value = "per_si_allsi.cov_q3_tot"
$ws.Cells.Item(17, 2).Value = "SI.COV.Q3.TOT"
$ws.Cells.Item(18, 1).Value = "per_si_allsi.cov_q4_tot"
$ws.Cells.Item(18, 2).Value = "SI.COV.Q4.TOT"
$ws.Cells.Item(19, 1).Value = "per_si_allsi.cov_q1_tot"
$ws.Cells.Item(19, 2).Value = "SI.COV.Q1.TOT"
$ws.Cells.Item(20, 1).Value = "per_sa_allsa.cov_q1_tot"
$ws.Cells.Item(20, 2).Value = "SA.COV.Q1.TOT"
$ws.Cells.Item(21, 1).Value = "per_si_allsi.cov_q5_tot"
$ws.Cells.Item(21, 2).Value = "SI.COV.Q5.TOT"
$ws.Cells.Item(22, 1).Value = "per_allsp.cov_pop_tot"
$ws.Cells.Item(22, 2).Value = "SP.COV.POP.TOT"
$ws.Cells.Item(23, 1).Value = "per_sa_allsa.cov_q3_tot"
$ws.Cells.Item(23, 2).Value = "SA.COV.Q3.TOT"
$ws.Cells.Item(24, 1).Value = "per_sa_allsa.cov_pop_tot"
$ws.Cells.Item(24, 2).Value = "SA.COV.POP.TOT"
$ws.Cells.Item(25, 1).Value = "per_sa_allsa.cov_q5_tot"
$ws.Cells.Item(25, 2).Value = "SA.COV.Q5.TOT"
$ws.Cells.Item(26, 1).Value = "per_lm_alllm.cov_q3_tot"
$ws.Cells.Item(26, 2).Value = "LM.COV.Q3.TOT"
$ws.Cells.Item(27, 1).Value = "per_lm_alllm.cov_q2_tot"
$ws.Cells.Item(27, 2).Value = "LM.COV.Q2.TOT"
$ws.Cells.Item(28, 1).Value = "per_lm_alllm.cov_pop_tot"
$ws.Cells.Item(28, 2).Value = "LM.COV.POP.TOT"
